# Update the NATMI ligand-receptor pair metrics for rows 2-17 following
# Dr Hou's advice (ligand/receptor expressing-cell counts go from 1 to 3
# "expression groups", which changes the downstream average/total
# expression, specificity and edge-weight values for every row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 122.253015
$ws.Range("H2").Value = 366.759045
$ws.Range("I2").Value = 0.1988639364328829
$ws.Range("J2").Value = 0.1988639364328829
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 57.478532
$ws.Range("N2").Value = 172.435596
$ws.Range("O2").Value = 0.2414676574042868
$ws.Range("P2").Value = 0.2414676574042868
$ws.Range("Q2").Value = 7026.92383477398
$ws.Range("R2").Value = 63242.31451296582
$ws.Range("S2").Value = 0.04801920887264323
$ws.Range("T2").Value = 0.04801920887264323

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 122.253015
$ws.Range("H3").Value = 366.759045
$ws.Range("I3").Value = 0.1988639364328829
$ws.Range("J3").Value = 0.1988639364328829
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 64.84043133333334
$ws.Range("N3").Value = 194.521294
$ws.Range("O3").Value = 0.272395040623924
$ws.Range("P3").Value = 0.2723950406239241
$ws.Range("Q3").Value = 7926.938224400471
$ws.Range("R3").Value = 71342.44401960424
$ws.Range("S3").Value = 0.05416955004326857
$ws.Range("T3").Value = 0.05416955004326858

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 122.253015
$ws.Range("H4").Value = 366.759045
$ws.Range("I4").Value = 0.1988639364328829
$ws.Range("J4").Value = 0.1988639364328829
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 85.31555666666667
$ws.Range("N4").Value = 255.94667
$ws.Range("O4").Value = 0.3584111648579104
$ws.Range("P4").Value = 0.3584111648579105
$ws.Range("Q4").Value = 10430.08402890335
$ws.Range("R4").Value = 93870.75626013015
$ws.Range("S4").Value = 0.07127505510513901
$ws.Range("T4").Value = 0.07127505510513901

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 122.253015
$ws.Range("H5").Value = 366.759045
$ws.Range("I5").Value = 0.1988639364328829
$ws.Range("J5").Value = 0.1988639364328829
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.40370266666666
$ws.Range("N5").Value = 91.211108
$ws.Range("O5").Value = 0.1277261371138787
$ws.Range("P5").Value = 0.1277261371138788
$ws.Range("Q5").Value = 3716.94431816354
$ws.Range("R5").Value = 33452.49886347186
$ws.Range("S5").Value = 0.02540012241183206
$ws.Range("T5").Value = 0.02540012241183207

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 132.5447616666667
$ws.Range("H6").Value = 397.634285
$ws.Range("I6").Value = 0.2156050961899926
$ws.Range("J6").Value = 0.2156050961899926
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 57.478532
$ws.Range("N6").Value = 172.435596
$ws.Range("O6").Value = 0.2414676574042868
$ws.Range("P6").Value = 0.2414676574042868
$ws.Range("Q6").Value = 7618.478324889873
$ws.Range("R6").Value = 68566.30492400886
$ws.Range("S6").Value = 0.05206165750142343
$ws.Range("T6").Value = 0.05206165750142343

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 132.5447616666667
$ws.Range("H7").Value = 397.634285
$ws.Range("I7").Value = 0.2156050961899926
$ws.Range("J7").Value = 0.2156050961899926
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 64.84043133333334
$ws.Range("N7").Value = 194.521294
$ws.Range("O7").Value = 0.272395040623924
$ws.Range("P7").Value = 0.2723950406239241
$ws.Range("Q7").Value = 8594.259517440532
$ws.Range("R7").Value = 77348.33565696479
$ws.Range("S7").Value = 0.05872975893539808
$ws.Range("T7").Value = 0.05872975893539809

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 132.5447616666667
$ws.Range("H8").Value = 397.634285
$ws.Range("I8").Value = 0.2156050961899926
$ws.Range("J8").Value = 0.2156050961899926
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 85.31555666666667
$ws.Range("N8").Value = 255.94667
$ws.Range("O8").Value = 0.3584111648579104
$ws.Range("P8").Value = 0.3584111648579105
$ws.Range("Q8").Value = 11308.13012484233
$ws.Range("R8").Value = 101773.1711235809
$ws.Range("S8").Value = 0.07727527367475708
$ws.Range("T8").Value = 0.07727527367475709

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 132.5447616666667
$ws.Range("H9").Value = 397.634285
$ws.Range("I9").Value = 0.2156050961899926
$ws.Range("J9").Value = 0.2156050961899926
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.40370266666666
$ws.Range("N9").Value = 91.211108
$ws.Range("O9").Value = 0.1277261371138787
$ws.Range("P9").Value = 0.1277261371138788
$ws.Range("Q9").Value = 4029.85152373753
$ws.Range("R9").Value = 36268.66371363778
$ws.Range("S9").Value = 0.02753840607841401
$ws.Range("T9").Value = 0.02753840607841401

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 320.0894206666666
$ws.Range("H10").Value = 960.2682619999999
$ws.Range("I10").Value = 0.5206762565675317
$ws.Range("J10").Value = 0.5206762565675317
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 57.478532
$ws.Range("N10").Value = 172.435596
$ws.Range("O10").Value = 0.2414676574042868
$ws.Range("P10").Value = 0.2414676574042868
$ws.Range("Q10").Value = 18398.27000865046
$ws.Range("R10").Value = 165584.4300778541
$ws.Range("S10").Value = 0.1257264759393953
$ws.Range("T10").Value = 0.1257264759393953

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 320.0894206666666
$ws.Range("H11").Value = 960.2682619999999
$ws.Range("I11").Value = 0.5206762565675317
$ws.Range("J11").Value = 0.5206762565675317
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 64.84043133333334
$ws.Range("N11").Value = 194.521294
$ws.Range("O11").Value = 0.272395040623924
$ws.Range("P11").Value = 0.2723950406239241
$ws.Range("Q11").Value = 20754.73610126345
$ws.Range("R11").Value = 186792.624911371
$ws.Range("S11").Value = 0.1418296300596255
$ws.Range("T11").Value = 0.1418296300596255

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 320.0894206666666
$ws.Range("H12").Value = 960.2682619999999
$ws.Range("I12").Value = 0.5206762565675317
$ws.Range("J12").Value = 0.5206762565675317
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 85.31555666666667
$ws.Range("N12").Value = 255.94667
$ws.Range("O12").Value = 0.3584111648579104
$ws.Range("P12").Value = 0.3584111648579105
$ws.Range("Q12").Value = 27308.6071072875
$ws.Range("R12").Value = 245777.4639655875
$ws.Range("S12").Value = 0.1866161836302253
$ws.Range("T12").Value = 0.1866161836302253

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 320.0894206666666
$ws.Range("H13").Value = 960.2682619999999
$ws.Range("I13").Value = 0.5206762565675317
$ws.Range("J13").Value = 0.5206762565675317
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.40370266666666
$ws.Range("N13").Value = 91.211108
$ws.Range("O13").Value = 0.1277261371138787
$ws.Range("P13").Value = 0.1277261371138788
$ws.Range("Q13").Value = 9731.903572694921
$ws.Range("R13").Value = 87587.13215425429
$ws.Range("S13").Value = 0.06650396693828566
$ws.Range("T13").Value = 0.06650396693828567

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 39.86989333333333
$ws.Range("H14").Value = 119.60968
$ws.Range("I14").Value = 0.06485471080959287
$ws.Range("J14").Value = 0.06485471080959287
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 57.478532
$ws.Range("N14").Value = 172.435596
$ws.Range("O14").Value = 0.2414676574042868
$ws.Range("P14").Value = 0.2414676574042868
$ws.Range("Q14").Value = 2291.662939796587
$ws.Range("R14").Value = 20624.96645816928
$ws.Range("S14").Value = 0.01566031509082487
$ws.Range("T14").Value = 0.01566031509082487

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 39.86989333333333
$ws.Range("H15").Value = 119.60968
$ws.Range("I15").Value = 0.06485471080959287
$ws.Range("J15").Value = 0.06485471080959287
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 64.84043133333334
$ws.Range("N15").Value = 194.521294
$ws.Range("O15").Value = 0.272395040623924
$ws.Range("P15").Value = 0.2723950406239241
$ws.Range("Q15").Value = 2585.181080947325
$ws.Range("R15").Value = 23266.62972852592
$ws.Range("S15").Value = 0.01766610158563189
$ws.Range("T15").Value = 0.0176661015856319

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 39.86989333333333
$ws.Range("H16").Value = 119.60968
$ws.Range("I16").Value = 0.06485471080959287
$ws.Range("J16").Value = 0.06485471080959287
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 85.31555666666667
$ws.Range("N16").Value = 255.94667
$ws.Range("O16").Value = 0.3584111648579104
$ws.Range("P16").Value = 0.3584111648579105
$ws.Range("Q16").Value = 3401.522143973955
$ws.Range("R16").Value = 30613.6992957656
$ws.Range("S16").Value = 0.0232446524477891
$ws.Range("T16").Value = 0.0232446524477891

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 39.86989333333333
$ws.Range("H17").Value = 119.60968
$ws.Range("I17").Value = 0.06485471080959287
$ws.Range("J17").Value = 0.06485471080959287
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.40370266666666
$ws.Range("N17").Value = 91.211108
$ws.Range("O17").Value = 0.1277261371138787
$ws.Range("P17").Value = 0.1277261371138788
$ws.Range("Q17").Value = 1212.192382258382
$ws.Range("R17").Value = 10909.73144032544
$ws.Range("S17").Value = 0.008283641685347013
$ws.Range("T17").Value = 0.008283641685347015
